$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = 30.4076393629437
$ws.Range("K2").Value = 1335.891102037011
$ws.Range("L2").Value = 36.54984407678111
$ws.Range("M2").Value = 0.1930345932205534
